# Update the multiplication equations in the document per the commit diff.
$d = $word.ActiveDocument

$replacements = @(
    @("89×51=4539", "36×66=2376"),
    @("59×26=1534", "23×84=1932"),
    @("59×29=1711", "56×83=4648"),
    @("90×26=2340", "60×53=3180"),
    @("66×19=1254", "40×43=1720"),
    @("90×75=6750", "42×77=3234"),
    @("48×56=2688", "26×55=1430"),
    @("55×92=5060", "51×86=4386"),
    @("42×68=2856", "22×86=1892"),
    @("97×23=2231", "19×54=1026"),
    @("99×12=1188", "25×34=850"),
    @("40×63=2520", "53×23=1219"),
    @("53×79=4187", "35×71=2485"),
    @("57×39=2223", "49×57=2793"),
    @("59×42=2478", "32×49=1568"),
    @("46×40=1840", "22×21=462"),
    @("80×51=4080", "53×51=2703"),
    @("30×88=2640", "41×67=2747"),
    @("27×11=297",  "60×31=1860"),
    @("93×98=9114", "44×58=2552"),
    @("67×60=4020", "17×63=1071"),
    @("34×86=2924", "17×76=1292"),
    @("22×28=616",  "90×84=7560"),
    @("24×26=624",  "79×99=7821"),
    @("77×61=4697", "15×37=555")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Replaced $($replacements.Count) equations."
